$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 5249.5
$ws.Range("I62").Value = 2999
$ws.Range("J62").Value = 7500
$ws.Range("K62").Value = 2999
$ws.Range("L62").Value = 7500
$ws.Range("M62").Value = -2375
$ws.Range("N62").Value = -8748

# Row 65
$ws.Range("H65").Value = 5249.5
$ws.Range("I65").Value = 2999
$ws.Range("J65").Value = 7500
$ws.Range("K65").Value = 14995
$ws.Range("L65").Value = 37500
$ws.Range("M65").Value = -11875
$ws.Range("N65").Value = -43740

# Row 137
$ws.Range("H137").Value = 5392.8213
$ws.Range("I137").Value = 1952.6086
$ws.Range("K137").Value = 5857.825800000001
$ws.Range("M137").Value = -3307.825800000001

$ws = $wb.Worksheets.Item("ARM")
# Row 14
$ws.Range("H14").Value = 2952
$ws.Range("J14").Value = 2952
$ws.Range("L14").Value = 2952
$ws.Range("N14").Value = -3302

# Row 16
$ws.Range("H16").Value = 131830.25
$ws.Range("I16").Value = 208245.6
$ws.Range("J16").Value = 4471.3335
$ws.Range("K16").Value = 208245.6
$ws.Range("L16").Value = 4471.3335
$ws.Range("M16").Value = -207958.6
$ws.Range("N16").Value = -5045.3335

# Row 61
$ws.Range("H61").Value = 12167.637
$ws.Range("I61").Value = 13983.223
$ws.Range("K61").Value = 13983.223
$ws.Range("M61").Value = -13771.223

# Row 104
$ws.Range("H104").Value = 7000
$ws.Range("J104").Value = 7000
$ws.Range("L104").Value = 7000
$ws.Range("N104").Value = -13988

# Row 122
$ws.Range("H122").Value = 2544.4
$ws.Range("I122").Value = 2115
$ws.Range("J122").Value = 2728.4285
$ws.Range("K122").Value = 6345
$ws.Range("L122").Value = 8185.2855
$ws.Range("M122").Value = -3895
$ws.Range("N122").Value = -13085.2855

# Row 132
$ws.Range("H132").Value = 4656.784
$ws.Range("I132").Value = 2665.6875
$ws.Range("K132").Value = 7997.0625
$ws.Range("M132").Value = -5467.0625

# Row 136
$ws.Range("H136").Value = 12167.637
$ws.Range("I136").Value = 13983.223
$ws.Range("K136").Value = 41949.669
$ws.Range("M136").Value = -39399.669

# Row 137
$ws.Range("H137").Value = 73964.664
$ws.Range("J137").Value = 73964.664
$ws.Range("L137").Value = 73964.664
$ws.Range("N137").Value = -84164.664

$ws = $wb.Worksheets.Item("BSM")
# Row 106
$ws.Range("H106").Value = 38333.332
$ws.Range("J106").Value = 38333.332
$ws.Range("L106").Value = 38333.332
$ws.Range("N106").Value = -40857.332

$ws = $wb.Worksheets.Item("CRP")
# Row 132
$ws.Range("H132").Value = 6694.6787
$ws.Range("I132").Value = 7819.1904
$ws.Range("J132").Value = 3321.1428
$ws.Range("K132").Value = 23457.5712
$ws.Range("L132").Value = 9963.428400000001
$ws.Range("M132").Value = -20927.5712
$ws.Range("N132").Value = -15023.4284

$ws = $wb.Worksheets.Item("CUL")
# Row 37
$ws.Range("H37").Value = 49614.617
$ws.Range("J37").Value = 49614.617
$ws.Range("L37").Value = 148843.851
$ws.Range("N37").Value = -149067.851

# Row 51
$ws.Range("H51").Value = 133336000
$ws.Range("J51").Value = 4000
$ws.Range("L51").Value = 12000
$ws.Range("N51").Value = -12920

# Row 132
$ws.Range("H132").Value = 4924.4287
$ws.Range("I132").Value = 2572.1667
$ws.Range("K132").Value = 23149.5003
$ws.Range("M132").Value = -20619.5003

$ws = $wb.Worksheets.Item("GSM")
# Row 18
$ws.Range("H18").Value = 10002500
$ws.Range("I18").Value = 10002500
$ws.Range("K18").Value = 10002500
$ws.Range("M18").Value = -10002207

# Row 21
$ws.Range("H21").Value = 4597.3335
$ws.Range("I21").Value = 4597.3335
$ws.Range("K21").Value = 4597.3335
$ws.Range("M21").Value = -4424.3335

# Row 30
$ws.Range("H30").Value = 4597.3335
$ws.Range("I30").Value = 4597.3335
$ws.Range("K30").Value = 4597.3335
$ws.Range("M30").Value = -4492.3335

# Row 32
$ws.Range("H32").Value = 29166.334
$ws.Range("J32").Value = 29166.334
$ws.Range("L32").Value = 29166.334
$ws.Range("N32").Value = -29758.334

# Row 42
$ws.Range("H42").Value = 57374.75
$ws.Range("J42").Value = 57374.75
$ws.Range("L42").Value = 57374.75
$ws.Range("N42").Value = -58344.75

# Row 115
$ws.Range("H115").Value = 57374.75
$ws.Range("J115").Value = 57374.75
$ws.Range("L115").Value = 57374.75
$ws.Range("N115").Value = -59724.75

# Row 122
$ws.Range("H122").Value = 2815.1428
$ws.Range("I122").Value = 2029.8889
$ws.Range("J122").Value = 3404.0833
$ws.Range("K122").Value = 6089.6667
$ws.Range("L122").Value = 10212.2499
$ws.Range("M122").Value = -3639.6667
$ws.Range("N122").Value = -15112.2499

# Row 132
$ws.Range("H132").Value = 10185.667
$ws.Range("I132").Value = 8512.941000000001
$ws.Range("J132").Value = 13029.3
$ws.Range("K132").Value = 25538.823
$ws.Range("L132").Value = 39087.89999999999
$ws.Range("M132").Value = -23008.823
$ws.Range("N132").Value = -44147.89999999999

$ws = $wb.Worksheets.Item("LTW")
# Row 25
$ws.Range("H25").Value = 6452.091
$ws.Range("I25").Value = 7068.625
$ws.Range("J25").Value = 4808
$ws.Range("K25").Value = 7068.625
$ws.Range("L25").Value = 4808
$ws.Range("M25").Value = -6838.625
$ws.Range("N25").Value = -5268

# Row 70
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

# Row 73
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

# Row 75
$ws.Range("H75").Value = 10851.667
$ws.Range("I75").Value = 10851.667
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 10851.667
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -9915.666999999999
$ws.Range("N75").ClearContents()

# Row 78
$ws.Range("H78").Value = 10851.667
$ws.Range("I78").Value = 10851.667
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 32555.001
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -27875.001
$ws.Range("N78").ClearContents()

# Row 106
$ws.Range("H106").Value = 8783.625
$ws.Range("J106").Value = 8783.625
$ws.Range("L106").Value = 8783.625
$ws.Range("N106").Value = -11307.625

# Row 122
$ws.Range("H122").Value = 8000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 24000
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -28900

$ws = $wb.Worksheets.Item("WVR")
# Row 23
$ws.Range("H23").Value = 1338.6666
$ws.Range("I23").Value = 259.8889
$ws.Range("J23").Value = 4575
$ws.Range("K23").Value = 259.8889
$ws.Range("L23").Value = 4575
$ws.Range("M23").Value = -30.88889999999998
$ws.Range("N23").Value = -5033

# Row 132
$ws.Range("H132").Value = 2236.9614
$ws.Range("I132").Value = 2126.48
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 6379.440000000001
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -3849.440000000001
$ws.Range("N132").Value = -20057

# Row 136
$ws.Range("H136").Value = 73857.57000000001
$ws.Range("J136").Value = 502777.5
$ws.Range("L136").Value = 1508332.5
$ws.Range("N136").Value = -1513432.5
